$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Sprite" column (K) - header + type row + per-row values
$ws.Range("K1").Value = "Sprite"
$ws.Range("K2").Value = "string"
$ws.Range("K3").Value = "Cookie/DummyCookie1"
$ws.Range("K4").Value = "Cookie/DummyCookie1"
$ws.Range("K5").Value = "Cookie/DummyCookie1"
$ws.Range("K6").Value = "Cookie/DummyCookie1"
$ws.Range("K7").Value = "Cookie/DummyCookie1"
$ws.Range("K8").Value = "Cookie/DummyCookie1"
$ws.Range("K9").Value = "Cookie/DummyCookie1"
$ws.Range("K10").Value = "Cookie/DummyCookie1"

# Update current selection to match the authored state
$ws.Range("L7").Select() | Out-Null
